$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "69.648.87"
Set-TextValue 2 5 "  -0.04%  "

# Row 3
Set-TextValue 3 4 "3.748.27"
Set-TextValue 3 5 "  +0.17%  "

# Row 4
Set-TextValue 4 5 "  +0.10%  "

# Row 5
Set-TextValue 5 4 "613.32"
Set-TextValue 5 5 "  +0.02%  "

# Row 6
Set-TextValue 6 4 "178.89"
Set-TextValue 6 5 "  +0.74%  "

# Row 7
Set-TextValue 7 4 "3.746.66"
Set-TextValue 7 5 "  +0.16%  "

# Row 8
Set-TextValue 8 5 "  -0.01%  "

# Row 9
Set-TextValue 9 4 "0.528"
Set-TextValue 9 5 "  -2.29%  "

# Row 10
Set-TextValue 10 5 "  -0.32%  "

# Row 11
Set-TextValue 11 5 "  +3.71%  "

# Row 12
Set-TextValue 12 4 "0.481"
Set-TextValue 12 5 "  -3.59%  "

# Row 13
Set-TextValue 13 4 "40.09"
Set-TextValue 13 5 "  -2.09%  "

# Row 14
Set-TextValue 14 4 "0.0000253"
Set-TextValue 14 5 "  -0.23%  "

# Row 15
Set-TextValue 15 4 "4.373.98"
Set-TextValue 15 5 "  +0.18%  "

# Row 16
Set-TextValue 16 4 "3.763.34"
Set-TextValue 16 5 "  +0.67%  "

# Row 17
Set-TextValue 17 4 "69.735.37"
Set-TextValue 17 5 "  -0.05%  "

# Row 19
Set-TextValue 19 4 "7.45"
Set-TextValue 19 5 "  -2.16%  "

# Row 20
Set-TextValue 20 4 "16.39"
Set-TextValue 20 5 "  -2.26%  "

# Row 21
Set-TextValue 21 4 "501.50"
Set-TextValue 21 5 "  -2.65%  "

# Row 22
Set-TextValue 22 4 "9.16"
Set-TextValue 22 5 "  -4.16%  "

# Row 23
Set-TextValue 23 4 "0.721"
Set-TextValue 23 5 "  -0.99%  "

# Row 24
Set-TextValue 24 4 "2.60"
Set-TextValue 24 5 "  +3.71%  "

# Row 25
Set-TextValue 25 4 "85.96"
Set-TextValue 25 5 "  -2.42%  "

# Row 26
Set-TextValue 26 4 "11.23"
Set-TextValue 26 5 "  +1.44%  "

# Row 27
Set-TextValue 27 4 "12.90"
Set-TextValue 27 5 "  -3.96%  "

# Row 28
Set-TextValue 28 4 "0.0000134"
Set-TextValue 28 5 "  +6.16%  "

# Row 29
Set-TextValue 29 5 "  -0.07%  "

# Row 30
Set-TextValue 30 2 "ImmutableX"
Set-TextValue 30 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 30 4 "2.47"
Set-TextValue 30 5 "  -1.65%  "

# Row 31
Set-TextValue 31 2 "NEARProtocol"
Set-TextValue 31 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 31 4 "8.05"
Set-TextValue 31 5 "  +2.32%  "

# Row 32
Set-TextValue 32 5 "  +2.31%  "

# Row 33
Set-TextValue 33 4 "30.41"
Set-TextValue 33 5 "  -2.73%  "

# Row 34
Set-TextValue 34 5 "  -2.12%  "

# Row 35
Set-TextValue 35 5 "  +0.16%  "

# Row 36
Set-TextValue 36 5 "  +1.17%  "

# Row 37
Set-TextValue 37 4 "6.11"
Set-TextValue 37 5 "  -1.93%  "

# Row 38
Set-TextValue 38 4 "0.350"
Set-TextValue 38 5 "  +2.75%  "

# Row 39
Set-TextValue 39 4 "0.137"
Set-TextValue 39 5 "  +3.07%  "

# Row 40
Set-TextValue 40 4 "449.12"
Set-TextValue 40 5 "  +6.20%  "

# Row 41
Set-TextValue 41 4 "3.06"
Set-TextValue 41 5 "  +12.55%  "

# Row 42
Set-TextValue 42 4 "2.06"
Set-TextValue 42 5 "  -5.31%  "

# Row 43
Set-TextValue 43 2 "OKB"
Set-TextValue 43 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 43 4 "49.71"
Set-TextValue 43 5 "  -3.20%  "

# Row 44
Set-TextValue 44 2 "Arweave"
Set-TextValue 44 3 "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue 44 4 "45.28"
Set-TextValue 44 5 "  +2.05%  "

# Row 45
Set-TextValue 45 4 "8.56"
Set-TextValue 45 5 "  -2.99%  "

# Row 46
Set-TextValue 46 4 "2.953.28"
Set-TextValue 46 5 "  -4.62%  "

# Row 47
Set-TextValue 47 4 "0.0359"
Set-TextValue 47 5 "  -1.68%  "

# Row 48
Set-TextValue 48 2 "Monero"
Set-TextValue 48 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 48 4 "139.11"
Set-TextValue 48 5 "  +2.82%  "

# Row 49
Set-TextValue 49 2 "USDe"
Set-TextValue 49 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 49 4 "1.00"
Set-TextValue 49 5 "  -0.06%  "

# Row 50
Set-TextValue 50 2 "InjectiveProtocol"
Set-TextValue 50 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 50 4 "27.09"
Set-TextValue 50 5 "  -3.00%  "

# Row 51
Set-TextValue 51 4 "2.49"
Set-TextValue 51 5 "  -1.83%  "
